$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of consultation data (row 3). Values that could otherwise be
# auto-coerced into numbers/dates/times by Excel's input parser (phone
# number, ISO date, time-of-day) are entered with a leading apostrophe so
# they are stored as literal text, matching the other rows in the sheet.
$ws.Range("A3").Value = "Nitigya"
$ws.Range("B3").Value = "'1897326235"
$ws.Range("C3").Value = "Laude me dard"
$ws.Range("D3").Value = "Dr. Linda Martinez: Gynecologist"
$ws.Range("E3").Value = "'2024-09-25"
$ws.Range("F3").Value = "'19:08"
